# Fix 3-Year Summary category alignment
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-Year Summary")

# Re-label existing category rows so that they line up with the order used
# on the Infrastructure Costs / Credits sheets:
#   row 3 -> Cloud Services
#   row 4 -> Software Licenses
#   row 5 -> Support & Maintenance
$ws.Range("A3").Value = "Cloud Services"
$ws.Range("A4").Value = "Software Licenses"
$ws.Range("A5").Value = "Support & Maintenance"

# Row 6 used to be the TOTAL row; turn it into a normal category row for the
# previously-missing "Professional Services" category. It already carries
# the right styling (s=51 / s=53), so just replace its contents.
$ws.Range("A6").Value = "Professional Services"
$ws.Range("B6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$G:`$G)"
$ws.Range("C6").Formula = "=SUMIF(Credits!`$A:`$A,A6,Credits!`$C:`$C)"
$ws.Range("D6").Formula = "=B6+C6"
$ws.Range("E6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$H:`$H)"
$ws.Range("F6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$I:`$I)"
$ws.Range("G6").Formula = "=D6+E6+F6"

# Insert the TOTAL row back in at row 7 (plain/default styling), summing the
# now four category rows (3-6).
$ws.Range("A7").Value = "TOTAL"

$ws.Range("B7").Formula = "=SUM(B3:B6)"
$ws.Range("C7").Formula = "=SUM(C3:C6)"
$ws.Range("D7").Formula = "=SUM(D3:D6)"
$ws.Range("E7").Formula = "=SUM(E3:E6)"
$ws.Range("F7").Formula = "=SUM(F3:F6)"
$ws.Range("G7").Formula = "=SUM(G3:G6)"

# Entering formulas that reference the styled currency cells above causes
# the engine to auto-copy their number format; strip it back off so row 7
# keeps the plain/default styling it has in the target workbook.
$ws.Range("B7:G7").ClearFormats()

# Row 8 stays blank but is now part of the sheet's used range (touching
# OutlineLevel materializes the row without adding any cells/attributes).
$ws.Rows.Item(8).OutlineLevel = 0
